# Perbaikan lebar tampilan profil + update data kendaraan dan user
# Append a new data row (row 4) to the user table on the active sheet:
#   NIK               | Plat | Nama        | Password
#   7868866666665555  |      | Siti Aminah | Akun_siti21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NIK looks numeric but must stay text (matches rows 2-3, which are stored
# as text), so force a text format before writing the value.
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "7868866666665555"

# Plat (B) stays blank for this user, same as rows 2-3 - copy the existing
# blank cell above so B4 keeps a real (empty) cell entry instead of being
# left completely absent from the sheet.
$ws.Range("B2").Copy($ws.Range("B4"))

$ws.Range("C4").Value = "Siti Aminah"
$ws.Range("D4").Value = "Akun_siti21"
